# edit.ps1 - apply KHL probabilities update (2025-11-24 -> 2025-11-25 slate)
# Generated to match the commit diff: updates Summary (3 match rows, all stat
# columns) and Cards_telegram (date/match/card_text) for rows 2-4.

$wb = $excel.ActiveWorkbook

# --- Sheet "Summary": update date, teams, ids, url and all probability/odds columns ---
$ws1 = $wb.Worksheets.Item("Summary")

# Row 2
$ws1.Cells.Item(2, 2).Value = 45986.79166666666
$ws1.Cells.Item(2, 3).Value = "Нефтехимик"
$ws1.Cells.Item(2, 4).Value = "Авангард"
$ws1.Cells.Item(2, 5).Value = "Нефтехимик – Авангард"
$ws1.Cells.Item(2, 6).Value = 897814
$ws1.Cells.Item(2, 7).Value = "https://text.khl.ru/text/897814.html"
$ws1.Cells.Item(2, 8).Value = 1.138889
$ws1.Cells.Item(2, 9).Value = 5.538462
$ws1.Cells.Item(2, 10).Value = 6.67735
$ws1.Cells.Item(2, 11).Value = 23.62224
$ws1.Cells.Item(2, 12).Value = 40.878477
$ws1.Cells.Item(2, 13).Value = 64.50071699999999
$ws1.Cells.Item(2, 14).Value = 0.040577
$ws1.Cells.Item(2, 15).Value = 0.059853
$ws1.Cells.Item(2, 16).Value = 0.890176
$ws1.Cells.Item(2, 17).Value = 24.64450304359612
$ws1.Cells.Item(2, 18).Value = 16.70760028737072
$ws1.Cells.Item(2, 19).Value = 1.123373355381408
$ws1.Cells.Item(2, 20).Value = 4.057700000000001
$ws1.Cells.Item(2, 21).Value = 5.985300000000001
$ws1.Cells.Item(2, 22).Value = 89.0176
$ws1.Cells.Item(2, 23).Value = 0.291673
$ws1.Cells.Item(2, 24).Value = 0.6989340000000001
$ws1.Cells.Item(2, 25).Value = 1.430750256819671
$ws1.Cells.Item(2, 26).Value = 0.453587
$ws1.Cells.Item(2, 27).Value = 0.5370200000000001
$ws1.Cells.Item(2, 28).Value = 1.862128039924025
$ws1.Cells.Item(2, 29).Value = 0.614178
$ws1.Cells.Item(2, 30).Value = 0.376429
$ws1.Cells.Item(2, 31).Value = 2.656543465035903
$ws1.Cells.Item(2, 32).Value = 0.345774
$ws1.Cells.Item(2, 33).Value = 0.654226
$ws1.Cells.Item(2, 34).Value = 2.892062445412321
$ws1.Cells.Item(2, 35).Value = 0.125602
$ws1.Cells.Item(2, 36).Value = 0.874398
$ws1.Cells.Item(2, 37).Value = 7.96165666151813
$ws1.Cells.Item(2, 38).Value = 0.949327
$ws1.Cells.Item(2, 39).Value = 0.050673
$ws1.Cells.Item(2, 40).Value = 1.053377813967158
$ws1.Cells.Item(2, 41).Value = 0.850458
$ws1.Cells.Item(2, 42).Value = 0.149542
$ws1.Cells.Item(2, 43).Value = 1.17583701958239
$ws1.Cells.Item(2, 44).Value = 0.205009
$ws1.Cells.Item(2, 45).Value = 4.877834631650318
$ws1.Cells.Item(2, 46).Value = 0.977089
$ws1.Cells.Item(2, 47).Value = 1.023448222219266

# Row 3
$ws1.Cells.Item(3, 2).Value = 45986.8125
$ws1.Cells.Item(3, 3).Value = "Спартак"
$ws1.Cells.Item(3, 4).Value = "Ак Барс"
$ws1.Cells.Item(3, 5).Value = "Спартак – Ак Барс"
$ws1.Cells.Item(3, 6).Value = 897812
$ws1.Cells.Item(3, 7).Value = "https://text.khl.ru/text/897812.html"
$ws1.Cells.Item(3, 8).Value = 4.588235
$ws1.Cells.Item(3, 9).Value = 3.665934
$ws1.Cells.Item(3, 10).Value = 8.254168999999999
$ws1.Cells.Item(3, 11).Value = 36.153363
$ws1.Cells.Item(3, 12).Value = 34.978634
$ws1.Cells.Item(3, 13).Value = 71.131997
$ws1.Cells.Item(3, 14).Value = 0.348824
$ws1.Cells.Item(3, 15).Value = 0.141178
$ws1.Cells.Item(3, 16).Value = 0.503965
$ws1.Cells.Item(3, 17).Value = 2.866775221888402
$ws1.Cells.Item(3, 18).Value = 7.083256598053521
$ws1.Cells.Item(3, 19).Value = 1.984264780292282
$ws1.Cells.Item(3, 20).Value = 34.8824
$ws1.Cells.Item(3, 21).Value = 14.1178
$ws1.Cells.Item(3, 22).Value = 50.3965
$ws1.Cells.Item(3, 23).Value = 0.10267
$ws1.Cells.Item(3, 24).Value = 0.891297
$ws1.Cells.Item(3, 25).Value = 1.121960468844841
$ws1.Cells.Item(3, 26).Value = 0.19608
$ws1.Cells.Item(3, 27).Value = 0.797887
$ws1.Cells.Item(3, 28).Value = 1.253310305845314
$ws1.Cells.Item(3, 29).Value = 0.319811
$ws1.Cells.Item(3, 30).Value = 0.674156
$ws1.Cells.Item(3, 31).Value = 1.483336201116656
$ws1.Cells.Item(3, 32).Value = 0.882715
$ws1.Cells.Item(3, 33).Value = 0.117285
$ws1.Cells.Item(3, 34).Value = 1.132868479633857
$ws1.Cells.Item(3, 35).Value = 0.712565
$ws1.Cells.Item(3, 36).Value = 0.287435
$ws1.Cells.Item(3, 37).Value = 1.403380744212809
$ws1.Cells.Item(3, 38).Value = 0.925683
$ws1.Cells.Item(3, 39).Value = 0.07431699999999999
$ws1.Cells.Item(3, 40).Value = 1.0802834231589
$ws1.Cells.Item(3, 41).Value = 0.797502
$ws1.Cells.Item(3, 42).Value = 0.202498
$ws1.Cells.Item(3, 43).Value = 1.253915350682506
$ws1.Cells.Item(3, 44).Value = 0.631793
$ws1.Cells.Item(3, 45).Value = 1.582796897084963
$ws1.Cells.Item(3, 46).Value = 0.767918
$ws1.Cells.Item(3, 47).Value = 1.302222372701252

# Row 4
$ws1.Cells.Item(4, 2).Value = 45986.8125
$ws1.Cells.Item(4, 3).Value = "ХК Сочи"
$ws1.Cells.Item(4, 4).Value = "Лада"
$ws1.Cells.Item(4, 5).Value = "ХК Сочи – Лада"
$ws1.Cells.Item(4, 6).Value = 897813
$ws1.Cells.Item(4, 7).Value = "https://text.khl.ru/text/897813.html"
$ws1.Cells.Item(4, 8).Value = 1.166667
$ws1.Cells.Item(4, 9).Value = 1.03125
$ws1.Cells.Item(4, 10).Value = 2.197917
$ws1.Cells.Item(4, 11).Value = 25.593223
$ws1.Cells.Item(4, 12).Value = 23.641989
$ws1.Cells.Item(4, 13).Value = 49.235212
$ws1.Cells.Item(4, 14).Value = 0.675169
$ws1.Cells.Item(4, 15).Value = 0.172768
$ws1.Cells.Item(4, 16).Value = 0.151996
$ws1.Cells.Item(4, 17).Value = 1.481110655258165
$ws1.Cells.Item(4, 18).Value = 5.788108909057232
$ws1.Cells.Item(4, 19).Value = 6.579120503171136
$ws1.Cells.Item(4, 20).Value = 67.51690000000001
$ws1.Cells.Item(4, 21).Value = 17.2768
$ws1.Cells.Item(4, 22).Value = 15.1996
$ws1.Cells.Item(4, 23).Value = 0.696346
$ws1.Cells.Item(4, 24).Value = 0.303588
$ws1.Cells.Item(4, 25).Value = 3.293937836805144
$ws1.Cells.Item(4, 26).Value = 0.836799
$ws1.Cells.Item(4, 27).Value = 0.163135
$ws1.Cells.Item(4, 28).Value = 6.129892420388022
$ws1.Cells.Item(4, 29).Value = 0.922303
$ws1.Cells.Item(4, 30).Value = 0.07763100000000001
$ws1.Cells.Item(4, 31).Value = 12.88145199726913
$ws1.Cells.Item(4, 32).Value = 0.716851
$ws1.Cells.Item(4, 33).Value = 0.283149
$ws1.Cells.Item(4, 34).Value = 1.394990032796216
$ws1.Cells.Item(4, 35).Value = 0.461394
$ws1.Cells.Item(4, 36).Value = 0.538606
$ws1.Cells.Item(4, 37).Value = 2.16734504566596
$ws1.Cells.Item(4, 38).Value = 0.312776
$ws1.Cells.Item(4, 39).Value = 0.6872239999999999
$ws1.Cells.Item(4, 40).Value = 3.197176253932527
$ws1.Cells.Item(4, 41).Value = 0.106172
$ws1.Cells.Item(4, 42).Value = 0.893828
$ws1.Cells.Item(4, 43).Value = 9.418679124439588
$ws1.Cells.Item(4, 44).Value = 0.944574
$ws1.Cells.Item(4, 45).Value = 1.058678303658581
$ws1.Cells.Item(4, 46).Value = 0.539866
$ws1.Cells.Item(4, 47).Value = 1.852311499520251

# --- Sheet "Cards_telegram": update date, match title and full card_text ---
$ws2 = $wb.Worksheets.Item("Cards_telegram")

# Row 2
$ws2.Cells.Item(2, 1).Value = 45986.79166666666
$ws2.Cells.Item(2, 2).Value = "Нефтехимик – Авангард"
$ws2.Cells.Item(2, 3).Value = "КХЛ • Регулярный чемпионат • 25.11.2025`n`nНефтехимик – Авангард`n`nОжидания модели (60’):`n• Голы: λ_total ≈ 5.95 (1.22 : 4.73)`n• Броски: SOG λ ≈ 65 (24 : 41)`n`nИсход (60’), честные кф:`n• П1: 4.1%  (Kмод 24.64)`n• Х:  6.0%  (Kмод 16.71)`n• П2: 89.0%  (Kмод 1.12)`n`nТоталы голов:`n• ТМ 4.5: 29.2%  (Kмод 3.43)`n• ТБ 4.5: 69.9%  (Kмод 1.43)`n`n• ТМ 5.5: 45.4%  (Kмод 2.20)`n• ТБ 5.5: 53.7%  (Kмод 1.86)`n`n• ТМ 6.5: 61.4%  (Kмод 1.63)`n• ТБ 6.5: 37.6%  (Kмод 2.66)`n`nИндивидуальные тоталы:`n• Нефтехимик ИТБ 1.5: 34.6% (Kмод 2.89)`n• Нефтехимик ИТБ 2.5: 12.6% (Kмод 7.96)`n• Авангард ИТБ 1.5: 94.9% (Kмод 1.05)`n• Авангард ИТБ 2.5: 85.0% (Kмод 1.18)`n`nФора +1.5:`n• Нефтехимик +1.5: 20.5% (Kмод 4.88)`n• Авангард +1.5: 97.7% (Kмод 1.02)"

# Row 3
$ws2.Cells.Item(3, 1).Value = 45986.8125
$ws2.Cells.Item(3, 2).Value = "Спартак – Ак Барс"
$ws2.Cells.Item(3, 3).Value = "КХЛ • Регулярный чемпионат • 25.11.2025`n`nСпартак – Ак Барс`n`nОжидания модели (60’):`n• Голы: λ_total ≈ 7.95 (3.69 : 4.26)`n• Броски: SOG λ ≈ 71 (36 : 35)`n`nИсход (60’), честные кф:`n• П1: 34.9%  (Kмод 2.87)`n• Х:  14.1%  (Kмод 7.08)`n• П2: 50.4%  (Kмод 1.98)`n`nТоталы голов:`n• ТМ 4.5: 10.3%  (Kмод 9.74)`n• ТБ 4.5: 89.1%  (Kмод 1.12)`n`n• ТМ 5.5: 19.6%  (Kмод 5.10)`n• ТБ 5.5: 79.8%  (Kмод 1.25)`n`n• ТМ 6.5: 32.0%  (Kмод 3.13)`n• ТБ 6.5: 67.4%  (Kмод 1.48)`n`nИндивидуальные тоталы:`n• Спартак ИТБ 1.5: 88.3% (Kмод 1.13)`n• Спартак ИТБ 2.5: 71.3% (Kмод 1.40)`n• Ак Барс ИТБ 1.5: 92.6% (Kмод 1.08)`n• Ак Барс ИТБ 2.5: 79.8% (Kмод 1.25)`n`nФора +1.5:`n• Спартак +1.5: 63.2% (Kмод 1.58)`n• Ак Барс +1.5: 76.8% (Kмод 1.30)"

# Row 4
$ws2.Cells.Item(4, 1).Value = 45986.8125
$ws2.Cells.Item(4, 2).Value = "ХК Сочи – Лада"
$ws2.Cells.Item(4, 3).Value = "КХЛ • Регулярный чемпионат • 25.11.2025`n`nХК Сочи – Лада`n`nОжидания модели (60’):`n• Голы: λ_total ≈ 3.65 (2.52 : 1.13)`n• Броски: SOG λ ≈ 49 (26 : 24)`n`nИсход (60’), честные кф:`n• П1: 67.5%  (Kмод 1.48)`n• Х:  17.3%  (Kмод 5.79)`n• П2: 15.2%  (Kмод 6.58)`n`nТоталы голов:`n• ТМ 4.5: 69.6%  (Kмод 1.44)`n• ТБ 4.5: 30.4%  (Kмод 3.29)`n`n• ТМ 5.5: 83.7%  (Kмод 1.20)`n• ТБ 5.5: 16.3%  (Kмод 6.13)`n`n• ТМ 6.5: 92.2%  (Kмод 1.08)`n• ТБ 6.5: 7.8%  (Kмод 12.88)`n`nИндивидуальные тоталы:`n• ХК Сочи ИТБ 1.5: 71.7% (Kмод 1.39)`n• ХК Сочи ИТБ 2.5: 46.1% (Kмод 2.17)`n• Лада ИТБ 1.5: 31.3% (Kмод 3.20)`n• Лада ИТБ 2.5: 10.6% (Kмод 9.42)`n`nФора +1.5:`n• ХК Сочи +1.5: 94.5% (Kмод 1.06)`n• Лада +1.5: 54.0% (Kмод 1.85)"

